$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 2251.75
$ws.Range("I62").Value2 = 2002.6666
$ws.Range("K62").Value2 = 2002.6666
$ws.Range("M62").Value2 = -1378.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 2251.75
$ws.Range("I65").Value2 = 2002.6666
$ws.Range("K65").Value2 = 10013.333
$ws.Range("M65").Value2 = -6893.333000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value2 = 5453.6665
$ws.Range("I70").Value2 = 5297.7144
$ws.Range("K70").Value2 = 15893.1432
$ws.Range("M70").Value2 = -15623.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value2 = 5453.6665
$ws.Range("I73").Value2 = 5297.7144
$ws.Range("K73").Value2 = 15893.1432
$ws.Range("M73").Value2 = -14957.1432

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value2 = 1354.421
$ws.Range("I80").Value2 = 867.44446
$ws.Range("J80").Value2 = 1792.7
$ws.Range("K80").Value2 = 2602.33338
$ws.Range("L80").Value2 = 5378.1
$ws.Range("M80").Value2 = -1604.33338
$ws.Range("N80").Value2 = -7374.1

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value2 = 1354.421
$ws.Range("I83").Value2 = 867.44446
$ws.Range("J83").Value2 = 1792.7
$ws.Range("K83").Value2 = 7807.00014
$ws.Range("L83").Value2 = 16134.3
$ws.Range("M83").Value2 = -2815.00014
$ws.Range("N83").Value2 = -26118.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value2 = 2392.5
$ws.Range("I98").Value2 = 2392.5
$ws.Range("K98").Value2 = 2392.5
$ws.Range("M98").Value2 = -894.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value2 = 2724.5
$ws.Range("I106").Value2 = 2632.6667
$ws.Range("K106").Value2 = 2632.6667
$ws.Range("M106").Value2 = -2001.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value2 = 2392.5
$ws.Range("I122").Value2 = 2392.5
$ws.Range("K122").Value2 = 7177.5
$ws.Range("M122").Value2 = -4727.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value2 = 10466.167
$ws.Range("I127").Value2 = 1400
$ws.Range("K127").Value2 = 4200
$ws.Range("M127").Value2 = 760

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value2 = 4909
$ws.Range("I138").Value2 = 1134.4
$ws.Range("J138").Value2 = 7370.696
$ws.Range("K138").Value2 = 3403.2
$ws.Range("L138").Value2 = 22112.088
$ws.Range("M138").Value2 = 1736.8
$ws.Range("N138").Value2 = -32392.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value2 = 201.14285
$ws.Range("I4").Value2 = 169.66667
$ws.Range("K4").Value2 = 169.66667
$ws.Range("M4").Value2 = -53.66667000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value2 = 3149.5
$ws.Range("I21").Value2 = 3149.5
$ws.Range("J21").Value2 = 0
$ws.Range("K21").Value2 = 3149.5
$ws.Range("L21").Value2 = 0
$ws.Range("M21").Value2 = -2775.5
$ws.Range("N21").Value2 = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value2 = 9500
$ws.Range("I23").Value2 = 8000
$ws.Range("K23").Value2 = 8000
$ws.Range("M23").Value2 = -7741

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value2 = 5574.75
$ws.Range("I102").Value2 = 4433
$ws.Range("K102").Value2 = 4433
$ws.Range("M102").Value2 = -2811

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 5336.5
$ws.Range("J99").Value2 = 5813.5
$ws.Range("L99").Value2 = 5813.5
$ws.Range("N99").Value2 = -8809.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 2793.4285
$ws.Range("I107").Value2 = 2793.4285
$ws.Range("K107").Value2 = 2793.4285
$ws.Range("M107").Value2 = -873.4285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 5816682
$ws.Range("I31").Value2 = 1938.5358
$ws.Range("J31").Value2 = 16670870
$ws.Range("K31").Value2 = 1938.5358
$ws.Range("L31").Value2 = 16670870
$ws.Range("M31").Value2 = -1643.5358
$ws.Range("N31").Value2 = -16671460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value2 = 28597.4
$ws.Range("I32").Value2 = 19000
$ws.Range("J32").Value2 = 34995.668
$ws.Range("K32").Value2 = 19000
$ws.Range("L32").Value2 = 34995.668
$ws.Range("M32").Value2 = -18684
$ws.Range("N32").Value2 = -35627.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 5816682
$ws.Range("I34").Value2 = 1938.5358
$ws.Range("J34").Value2 = 16670870
$ws.Range("K34").Value2 = 1938.5358
$ws.Range("L34").Value2 = 16670870
$ws.Range("M34").Value2 = -1736.5358
$ws.Range("N34").Value2 = -16671274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value2 = 3646.9524
$ws.Range("J134").Value2 = 4777.4287
$ws.Range("L134").Value2 = 14332.2861
$ws.Range("N134").Value2 = -19402.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value2 = 85680.65
$ws.Range("J141").Value2 = 88611.21
$ws.Range("L141").Value2 = 88611.21
$ws.Range("N141").Value2 = -98971.21

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value2 = 500.5
$ws.Range("I19").Value2 = 501
$ws.Range("J19").Value2 = 500
$ws.Range("K19").Value2 = 1503
$ws.Range("L19").Value2 = 1500
$ws.Range("M19").Value2 = -1329
$ws.Range("N19").Value2 = -1848

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value2 = 1000000
$ws.Range("I32").Value2 = 1000000
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 3000000
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -2999717
$ws.Range("N32").Value2 = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value2 = 1378.5714
$ws.Range("I39").Value2 = 895.93335
$ws.Range("J39").Value2 = 2585.1667
$ws.Range("K39").Value2 = 2687.80005
$ws.Range("L39").Value2 = 7755.500100000001
$ws.Range("M39").Value2 = -2393.80005
$ws.Range("N39").Value2 = -8343.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value2 = 3717.625
$ws.Range("I109").Value2 = 5926.857
$ws.Range("J109").Value2 = 1999.3334
$ws.Range("K109").Value2 = 17780.571
$ws.Range("L109").Value2 = 5998.0002
$ws.Range("M109").Value2 = -16740.571
$ws.Range("N109").Value2 = -8078.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value2 = 66669810
$ws.Range("I124").Value2 = 3924.75
$ws.Range("J124").Value2 = 333333340
$ws.Range("K124").Value2 = 11774.25
$ws.Range("L124").Value2 = 1000000020
$ws.Range("M124").Value2 = -6864.25
$ws.Range("N124").Value2 = -1000009840

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value2 = 69454360
$ws.Range("I126").Value2 = 125004130
$ws.Range("J126").Value2 = 41679480
$ws.Range("K126").Value2 = 375012390
$ws.Range("L126").Value2 = 125038440
$ws.Range("M126").Value2 = -375007450
$ws.Range("N126").Value2 = -125048320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value2 = 179995
$ws.Range("I128").Value2 = 179995
$ws.Range("K128").Value2 = 539985
$ws.Range("M128").Value2 = -535005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value2 = 2500
$ws.Range("I31").Value2 = 2500
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 2500
$ws.Range("L31").Value2 = 0
$ws.Range("M31").Value2 = -2208
$ws.Range("N31").Value2 = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value2 = 2500
$ws.Range("I37").Value2 = 2500
$ws.Range("J37").Value2 = 0
$ws.Range("K37").Value2 = 2500
$ws.Range("L37").Value2 = 0
$ws.Range("M37").Value2 = -2223
$ws.Range("N37").Value2 = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 8918.167
$ws.Range("I102").Value2 = 1701.8
$ws.Range("K102").Value2 = 1701.8
$ws.Range("M102").Value2 = -79.79999999999995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 38462736
$ws.Range("I122").Value2 = 1512.4286
$ws.Range("K122").Value2 = 4537.2858
$ws.Range("M122").Value2 = -2087.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 1200
$ws.Range("I132").Value2 = 1200
$ws.Range("K132").Value2 = 3600
$ws.Range("M132").Value2 = -1070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value2 = 20521.666
$ws.Range("J43").Value2 = 21282.5
$ws.Range("L43").Value2 = 21282.5
$ws.Range("N43").Value2 = -21668.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value2 = 3008.4878
$ws.Range("I100").Value2 = 2790.7812
$ws.Range("J100").Value2 = 3782.5557
$ws.Range("K100").Value2 = 2790.7812
$ws.Range("L100").Value2 = 3782.5557
$ws.Range("M100").Value2 = -2249.7812
$ws.Range("N100").Value2 = -4864.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value2 = 1922.3928
$ws.Range("I132").Value2 = 1737.4
$ws.Range("K132").Value2 = 5212.200000000001
$ws.Range("M132").Value2 = -2682.200000000001
